$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 1563
$ws1.Range("F14").Value = 5158
$ws1.Range("F17").Value = 229
$ws1.Range("F38").Value = 3

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 1563
$ws4.Range("F15").Value = 5158
$ws4.Range("F18").Value = 229
$ws4.Range("F39").Value = 3
